# Regenerate save_data: update column G ("K", formerly Strike#) values
# for the gillaspie_logan sheet. Rows 2-19 on Sheet1, column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    7  = 2
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
